# Update "想去人数" (F column) counts for both the "展览" and "全部类型" sheets.
# These two sheets contain identical data tables, and both need the same updates.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    2  = 2040
    7  = 1691
    9  = 688
    11 = 104
    12 = 27
    13 = 101
    16 = 140
    19 = 3919
    21 = 24
    23 = 368
    24 = 830
    25 = 566
    27 = 35
    28 = 1719
    31 = 174
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
